$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 494
$ws.Range("A3").Value = 497
$ws.Range("A4").Value = 499
$ws.Range("A5").Value = 500
$ws.Range("A6").Value = 502
$ws.Range("A7").Value = 503
$ws.Range("A8").Value = 505
$ws.Range("A9").Value = 509
$ws.Range("A10").Value = 511
$ws.Range("A11").Value = 513
$ws.Range("A12").Value = 514
$ws.Range("A13").Value = 516
$ws.Range("A14").Value = 517
$ws.Range("A15").Value = 7
$ws.Range("A16").Value = 30
$ws.Range("A17").Value = 62
$ws.Range("A18").Value = 119
$ws.Range("A19").Value = 155
$ws.Range("A20").Value = 175
$ws.Range("A21").Value = 220
$ws.Range("A22").Value = 235
$ws.Range("A23").Value = 282
$ws.Range("A24").Value = 336
$ws.Range("A25").Value = 366
$ws.Range("A26").Value = 407
$ws.Range("A27").Value = 481
